# Updated symbol list on Mon Jan 16 20:36:30 UTC 2023 with GitHub Actions
#
# Refresh the Price (D) and Volume(1h) (E) columns with the latest scraped
# quotes. Every value on this sheet is stored as literal text (inline
# strings, e.g. "299.90" / "-0.29%"), so each new value is entered with a
# leading apostrophe to force a text literal instead of letting Excel
# auto-convert the numeric-looking text into a real number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    if ($null -ne $value) {
        $ws.Range($cell).Value = "'" + $value
    }
}

# row, Price(D), Volume1h(E)
$updates = @(
    @(2,  "300.32",        "-0.13%"),
    @(3,  "31.73",         "1.19%"),
    @(4,  "5.140",         "0.56%"),
    @(5,  "0.08213",       "11.60%"),
    @(6,  "2.557",         "8.87%"),
    @(7,  "7.857",         "-1.18%"),
    @(8,  "3.855",         "1.77%"),
    @(9,  "0.9284",        "1.12%"),
    @(10, $null,           "3.19%"),
    @(11, "0.07510",       "-0.94%"),
    @(12, "0.09050",       "11.69%"),
    @(13, "0.03017",       "-0.05%"),
    @(14, $null,           "0.89%"),
    @(15, $null,           "2.01%"),
    @(16, "0.005994",      "-1.58%"),
    @(17, "3.620",         "4.49%"),
    @(18, $null,           "2.61%"),
    @(19, "0.3262",        "-0.67%"),
    @(20, "0.1346",        "0.70%"),
    @(21, "4.252",         "-8.59%"),
    @(22, "0.1680",        "7.24%"),
    @(23, "0.04631",       "-0.02%"),
    @(24, "0.001245",      "1.56%"),
    @(25, "0.004561",      "1.86%"),
    @(26, "0.0001198",     "-7.78%"),
    @(27, "0.0003399",     "81.61%"),
    @(39, $null,           "3.70%"),
    @(40, "0.04615",       "2.13%"),
    @(41, "0.006892",      "-4.22%"),
    @(42, $null,           "2.61%"),
    @(43, "0.002136",      "-4.56%"),
    @(44, "0.009828",      "-8.37%"),
    @(45, "0.00006166",    "-1.69%"),
    @(46, "0.00000000748", "-0.33%"),
    @(47, "0.8036",        "-0.61%"),
    @(48, "0.008373",      "-16.20%"),
    @(49, "0.00002093",    "-0.33%"),
    @(50, "0.0001994",     "-0.26%")
)

foreach ($u in $updates) {
    $row = $u[0]
    $price = $u[1]
    $volume = $u[2]

    Set-TextValue "D$row" $price
    Set-TextValue "E$row" $volume
}
